$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 2, 1.489477579669142),
    @(2, 3, 0.1874954686519459),
    @(2, 5, 0.1772940359760735),
    @(2, 6, 2.033009967590601),
    @(2, 7, 0.976411116029638),
    @(2, 8, 0.978897974828385),
    @(2, 9, 0.9170943720519134),
    @(2, 10, 0.0405182314035617),
    @(2, 12, 0.511870541956128),
    @(2, 14, 1.321400642053796),
    @(3, 2, 1.382237922499371),
    @(3, 3, 0.1690178688503465),
    @(3, 5, 0.1776175771092277),
    @(3, 6, 2.025088319577449),
    @(3, 7, 0.9684340870399524),
    @(3, 8, 0.9808816365426622),
    @(3, 9, 0.9220049875767771),
    @(3, 10, 0.04072690103194709),
    @(3, 12, 0.5018974493556527),
    @(3, 14, 1.337022150083595),
    @(4, 2, 1.31685600328683),
    @(4, 3, 0.1575937255937845),
    @(4, 5, 0.1778552994635447),
    @(4, 6, 2.021382673448798),
    @(4, 7, 0.9642910828353166),
    @(4, 8, 0.9826548267875523),
    @(4, 9, 0.9255878345730011),
    @(4, 10, 0.04086577053323737),
    @(4, 12, 0.4959734354524272),
    @(4, 14, 1.347155527685434),
    @(5, 2, 1.290329938994091),
    @(5, 3, 0.1529184232514069),
    @(5, 5, 0.177962026972196),
    @(5, 6, 2.020163688026514),
    @(5, 7, 0.9627920279541371),
    @(5, 8, 0.9835168870206843),
    @(5, 9, 0.9271904005520568),
    @(5, 10, 0.04092506984614452),
    @(5, 12, 0.4936097077697497),
    @(5, 14, 1.351421148911022),
    @(6, 2, 1.285932438778445),
    @(6, 3, 0.1521408912486493),
    @(6, 5, 0.177980345019737),
    @(6, 6, 2.019978853204648),
    @(6, 7, 0.9625545245609572),
    @(6, 8, 0.9836684510394917),
    @(6, 9, 0.9274651070212698),
    @(6, 10, 0.040935080262269),
    @(6, 12, 0.4932202595555566),
    @(6, 14, 1.35213767669579),
    @(7, 2, 1.316497786017123),
    @(7, 3, 0.1575307533524892),
    @(7, 5, 0.177856698890329),
    @(7, 6, 2.02136505537014),
    @(7, 7, 0.9642701005477932),
    @(7, 8, 0.9826658883140738),
    @(7, 9, 0.9256088705590315),
    @(7, 10, 0.04086655928807126),
    @(7, 12, 0.495941353251709),
    @(7, 14, 1.347212503967366),
    @(8, 2, 1.452405678572006),
    @(8, 3, 0.1811407159188718),
    @(8, 5, 0.177397502521135),
    @(8, 6, 2.03003807582742),
    @(8, 7, 0.9735035977716819),
    @(8, 8, 0.9794666185248673),
    @(8, 9, 0.9186695852987441),
    @(8, 10, 0.04058795452127129),
    @(8, 12, 0.5083905212026281),
    @(8, 14, 1.326674488912662),
    @(9, 2, 1.722570100097585),
    @(9, 3, 0.2268187768092673),
    @(9, 5, 0.1768057310215383),
    @(9, 6, 2.05624747587774),
    @(9, 7, 0.9976290615348944),
    @(9, 8, 0.9776057630002128),
    @(9, 9, 0.9095771744597982),
    @(9, 10, 0.04012657718739909),
    @(9, 12, 0.5343797243926218),
    @(9, 14, 1.290698710665545),
    @(10, 2, 1.923266803708714),
    @(10, 3, 0.2600090460544777),
    @(10, 5, 0.176557532623562),
    @(10, 6, 2.081134208815769),
    @(10, 7, 1.01906638807759),
    @(10, 8, 0.9789408500029992),
    @(10, 9, 0.9056659610250151),
    @(10, 10, 0.03983901075053886),
    @(10, 12, 0.5544287654267066),
    @(10, 14, 1.266889112425872),
    @(11, 2, 2.015045827427343),
    @(11, 3, 0.2750301118677498),
    @(11, 5, 0.1764847977807449),
    @(11, 6, 2.093683434325342),
    @(11, 7, 1.029634449953363),
    @(11, 8, 0.9801377190178187),
    @(11, 9, 0.9044914938363249),
    @(11, 10, 0.03971927194413638),
    @(11, 12, 0.5637557641568947),
    @(11, 14, 1.256627092022406),
    @(12, 2, 2.049868666452767),
    @(12, 3, 0.2807071810991033),
    @(12, 5, 0.1764630029410519),
    @(12, 6, 2.09861240327551),
    @(12, 7, 1.033754363330985),
    @(12, 8, 0.9806759155527232),
    @(12, 9, 0.9041339989424984),
    @(12, 10, 0.039675516532375),
    @(12, 12, 0.567317213837029),
    @(12, 14, 1.252823027644048),
    @(13, 2, 2.042365932453322),
    @(13, 3, 0.2794850148929982),
    @(13, 5, 0.176467441644574),
    @(13, 6, 2.09754299223043),
    @(13, 7, 1.032861806842277),
    @(13, 8, 0.980556223196686),
    @(13, 9, 0.9042071068073199),
    @(13, 10, 0.03968486954741657),
    @(13, 12, 0.5665488821272362),
    @(13, 14, 1.253638654911221),
    @(14, 2, 2.017909364302398),
    @(14, 3, 0.2754973896392414),
    @(14, 5, 0.1764828896644239),
    @(14, 6, 2.09408539781586),
    @(14, 7, 1.029971028658707),
    @(14, 8, 0.9801802927609629),
    @(14, 9, 0.9044603322450442),
    @(14, 10, 0.03971564038494257),
    @(14, 12, 0.5640481762693099),
    @(14, 14, 1.256312486527484),
    @(15, 2, 2.002937849218142),
    @(15, 3, 0.2730534120115635),
    @(15, 5, 0.1764930997767458),
    @(15, 6, 2.091990559804472),
    @(15, 7, 1.028215732611642),
    @(15, 8, 0.9799610956862068),
    @(15, 9, 0.9046268112263718),
    @(15, 10, 0.03973469492795445),
    @(15, 12, 0.5625202601896433),
    @(15, 14, 1.257960959860768),
    @(16, 2, 1.917278418106321),
    @(16, 3, 0.2590258310000877),
    @(16, 5, 0.1765630919457397),
    @(16, 6, 2.080338821465489),
    @(16, 7, 1.018392219177002),
    @(16, 8, 0.9788745112390984),
    @(16, 9, 0.9057549090221784),
    @(16, 10, 0.03984705836021973),
    @(16, 12, 0.5538233649971716),
    @(16, 14, 1.267571219994711),
    @(17, 2, 1.864851588783267),
    @(17, 3, 0.25040058103761),
    @(17, 5, 0.1766162985894493),
    @(17, 6, 2.073505584502712),
    @(17, 7, 1.012575308251158),
    @(17, 8, 0.9783590576017502),
    @(17, 9, 0.9066020579460812),
    @(17, 10, 0.03991882255561485),
    @(17, 12, 0.5485408787271666),
    @(17, 14, 1.273612633471597),
    @(18, 2, 1.834742408431907),
    @(18, 3, 0.2454322819447441),
    @(18, 5, 0.1766506844970532),
    @(18, 6, 2.06969086942523),
    @(18, 7, 1.009306351800575),
    @(18, 8, 0.9781180600073469),
    @(18, 9, 0.9071462215742088),
    @(18, 10, 0.03996114233867587),
    @(18, 12, 0.5455219887859499),
    @(18, 14, 1.277141052969995),
    @(19, 2, 1.824555778005163),
    @(19, 3, 0.2437488490630813),
    @(19, 5, 0.1766629775508122),
    @(19, 6, 2.06841911631976),
    @(19, 7, 1.008212705704523),
    @(19, 8, 0.9780459848313114),
    @(19, 9, 0.907340231000056),
    @(19, 10, 0.03997565038854844),
    @(19, 12, 0.5445031919808798),
    @(19, 14, 1.278344910731814),
    @(20, 2, 1.870427828944173),
    @(20, 3, 0.2513195071373673),
    @(20, 5, 0.1766102433469765),
    @(20, 6, 2.074221030147228),
    @(20, 7, 1.013186577445552),
    @(20, 8, 0.9784081853806015),
    @(20, 9, 0.9065059858707158),
    @(20, 10, 0.03991107523825121),
    @(20, 12, 0.5491011961178742),
    @(20, 14, 1.272963970857823),
    @(21, 2, 2.025091011036181),
    @(21, 3, 0.2766689522791239),
    @(21, 5, 0.1764781964255455),
    @(21, 6, 2.095096175482155),
    @(21, 7, 1.030816912034908),
    @(21, 8, 0.9802884050705813),
    @(21, 9, 0.9043835835578662),
    @(21, 10, 0.03970655921632549),
    @(21, 12, 0.5647818946019925),
    @(21, 14, 1.255524892662088),
    @(22, 2, 2.126568877562079),
    @(22, 3, 0.2931717031657968),
    @(22, 5, 0.1764253883712641),
    @(22, 6, 2.109770220958566),
    @(22, 7, 1.043027561034535),
    @(22, 8, 0.9820125887829363),
    @(22, 9, 0.9035051666878502),
    @(22, 10, 0.03958214430004858),
    @(22, 12, 0.5752021019608691),
    @(22, 14, 1.244605109583183),
    @(23, 2, 2.07237229684705),
    @(23, 3, 0.2843697694052594),
    @(23, 5, 0.1764505178184343),
    @(23, 6, 2.10184399425448),
    @(23, 7, 1.036447317780215),
    @(23, 8, 0.9810469700500732),
    @(23, 9, 0.9039273536490029),
    @(23, 10, 0.03964770253670125),
    @(23, 12, 0.569624968795452),
    @(23, 14, 1.250389460554928),
    @(24, 2, 1.867906711935177),
    @(24, 3, 0.2509040903348989),
    @(24, 5, 0.1766129690945384),
    @(24, 6, 2.073897222690945),
    @(24, 7, 1.01290998824868),
    @(24, 8, 0.9783858023253913),
    @(24, 9, 0.9065492421251804),
    @(24, 10, 0.03991457449257929),
    @(24, 12, 0.5488478202909732),
    @(24, 14, 1.273257059496078),
    @(25, 2, 1.649094651751682),
    @(25, 3, 0.2145271886062687),
    @(25, 5, 0.1769329396010733),
    @(25, 6, 2.048169867678354),
    @(25, 7, 0.9904536902758565),
    @(25, 8, 0.9776355060797215),
    @(25, 9, 0.911551827640821),
    @(25, 10, 0.04024233660680387),
    @(25, 12, 0.5271808402993798),
    @(25, 14, 1.299970782457617),
)

foreach ($item in $data) {
    $ws.Cells.Item($item[0], $item[1]).Value = $item[2]
}
